# Add payment 71652621 (Cash) 2025-08-15T09:47:50
#
# The "phone" column (A) normally stores the phone number as a genuine
# number (see rows 2-5). Row 6 was an exception -- its phone number had been
# entered as text. This edit:
#   1. Normalizes A6 back to a real number (71652621), matching the rest of
#      the column.
#   2. Appends a new row 7 for the latest payment from that same phone
#      number. The new row keeps the phone number as text (the same way row
#      6 originally held it) and records amount/method/timestamp as usual.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: convert A6 from text back to a plain number -------------------
$ws.Range("A6").Value = 71652621

# --- Row 7: the new payment record ----------------------------------------
# Leading apostrophe forces Excel to store the digits as text rather than
# auto-converting the numeric-looking string to a number.
$ws.Range("A7").Value = "'71652621"
$ws.Range("B7").Value = 20
$ws.Range("C7").Value = "Cash"
$ws.Range("D7").Value = "2025-08-15T09:47:50"
